$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting rows 192:223 down to 193:224
$ws.Rows.Item(192).Insert()

# Fill in the new row 192 with the new data record
$ws.Cells.Item(192, 1).Value = 10
$ws.Cells.Item(192, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(192, 3).Value = "La Araucanía"
$ws.Cells.Item(192, 4).Value = 45127
$ws.Cells.Item(192, 4).NumberFormat = $ws.Cells.Item(193, 4).NumberFormat
$ws.Cells.Item(192, 5).Value = 9
$ws.Cells.Item(192, 6).Value = 100114002
$ws.Cells.Item(192, 7).Value = "Camote"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 80
$ws.Cells.Item(192, 11).Value = 26000
$ws.Cells.Item(192, 12).Value = 26000
$ws.Cells.Item(192, 13).Value = 26000
$ws.Cells.Item(192, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(192, 15).Value = "Perú"
$ws.Cells.Item(192, 16).Value = 1444
$ws.Cells.Item(192, 17).Value = 18
$ws.Cells.Item(192, 18).Value = "Hortaliza"
